$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "remove 803 from andrea" — row 94 is Andreas' row with D94 = 803 and an
# empty C94. Delete the entire row; rows 95:173 shift up to 94:172.
$ws.Rows.Item(94).Delete()

# Re-establish the autofilter over the new (one-row-smaller) data range.
# Turning AutoFilterMode off first avoids the no-op "toggle" semantics of
# calling .AutoFilter() on a range that already has an active filter.
$ws.AutoFilterMode = $false
$ws.Range("A1:D172").AutoFilter(1)

# The hidden _FilterDatabase defined name also needs to track the shrunk range.
$filterDatabaseName = $wb.Names.Item(1)
$filterDatabaseName.RefersTo = "=Sheet1!`$A`$1:`$D`$172"

# Mirror Excel's default post-row-delete selection: the whole row that slid
# up into the deleted row's place (now row 94) becomes selected.
$ws.Range("A94:XFD94").Select()
